# Scheduled-runner market-data refresh for the Gungnir_Profits leve-profit tracker.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per leve row
# across all eight crafting-job sheets with freshly pulled Universalis prices.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 8: On the Drip / Eye Drops
$ws.Range("H8").Value = 119
$ws.Range("I8").Value = 119
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 357
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -218
$ws.Range("N8").ClearContents()
# Row 18: You Grow, Girl / Growth Formula Beta
$ws.Range("H18").Value = 297.18182
$ws.Range("I18").Value = 316.55554
$ws.Range("J18").Value = 210
$ws.Range("K18").Value = 316.55554
$ws.Range("L18").Value = 210
$ws.Range("M18").Value = -32.55554000000001
$ws.Range("N18").Value = -778
# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 3270.7917
$ws.Range("I64").Value = 3208.25
$ws.Range("K64").Value = 3208.25
$ws.Range("M64").Value = -2960.25
# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 3270.7917
$ws.Range("I67").Value = 3208.25
$ws.Range("K67").Value = 3208.25
$ws.Range("M67").Value = -2350.25
# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 17552134
$ws.Range("I76").Value = 10235.214
$ws.Range("J76").Value = 66669452
$ws.Range("K76").Value = 10235.214
$ws.Range("L76").Value = 66669452
$ws.Range("M76").Value = -9920.214
$ws.Range("N76").Value = -66670082
# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 17552134
$ws.Range("I79").Value = 10235.214
$ws.Range("J79").Value = 66669452
$ws.Range("K79").Value = 10235.214
$ws.Range("L79").Value = 66669452
$ws.Range("M79").Value = -9143.214
$ws.Range("N79").Value = -66671636
# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 14: Waste Not, Want Not / Bronze Celata
$ws.Range("H14").Value = 161
$ws.Range("I14").Value = 161
$ws.Range("K14").Value = 161
$ws.Range("M14").Value = 14
# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 1828.4
$ws.Range("I63").Value = 1764.8889
$ws.Range("K63").Value = 1764.8889
$ws.Range("M63").Value = -1078.8889
# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 1828.4
$ws.Range("I66").Value = 1764.8889
$ws.Range("K66").Value = 8824.4445
$ws.Range("M66").Value = -5392.4445
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 21740942
$ws.Range("I74").Value = 22729012
$ws.Range("J74").Value = 3400
$ws.Range("K74").Value = 22729012
$ws.Range("L74").Value = 3400
$ws.Range("M74").Value = -22728138
$ws.Range("N74").Value = -5148
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 21740942
$ws.Range("I77").Value = 22729012
$ws.Range("J77").Value = 3400
$ws.Range("K77").Value = 113645060
$ws.Range("L77").Value = 17000
$ws.Range("M77").Value = -113640692
$ws.Range("N77").Value = -25736
# Row 88: The Mast Chance / Adamantite Rivets
$ws.Range("H88").Value = 2515.5625
$ws.Range("I88").Value = 2593.75
$ws.Range("J88").Value = 2437.375
$ws.Range("K88").Value = 2593.75
$ws.Range("L88").Value = 2437.375
$ws.Range("M88").Value = -2187.75
$ws.Range("N88").Value = -3249.375
# Row 91: The Rose and the Riveter (L) / Adamantite Rivets
$ws.Range("H91").Value = 2515.5625
$ws.Range("I91").Value = 2593.75
$ws.Range("J91").Value = 2437.375
$ws.Range("K91").Value = 2593.75
$ws.Range("L91").Value = 2437.375
$ws.Range("M91").Value = -1189.75
$ws.Range("N91").Value = -5245.375
# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 1123.1765
$ws.Range("I122").Value = 941.1667
$ws.Range("J122").Value = 1560
$ws.Range("K122").Value = 2823.5001
$ws.Range("L122").Value = 4680
$ws.Range("M122").Value = -373.5001000000002
$ws.Range("N122").Value = -9580
# Row 129: In-kweh-dible Cooking / Manganese Chocobo Frypan
$ws.Range("H129").Value = 25040.375
$ws.Range("I129").Value = 10909
$ws.Range("J129").Value = 48592.668
$ws.Range("K129").Value = 10909
$ws.Range("L129").Value = 48592.668
$ws.Range("M129").Value = -5909
$ws.Range("N129").Value = -58592.668
# Row 131: Additions to the Armoire / Chondrite Top of Maiming
$ws.Range("H131").Value = 120715
$ws.Range("J131").Value = 120715
$ws.Range("L131").Value = 120715
$ws.Range("N131").Value = -130795
# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 16: Port of Call: Ul'dah / Bronze Knuckles
$ws.Range("H16").Value = 40000
$ws.Range("J16").Value = 40000
$ws.Range("L16").Value = 40000
$ws.Range("N16").Value = -40340
# Row 135: Axes to the Maxes / Ruthenium War Axe
$ws.Range("H135").Value = 48355.555
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 48355.555
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 48355.555
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -58495.555
# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 3: Touch and Heal / Maple Pattens
$ws.Range("H3").Value = 1001.5
$ws.Range("J3").Value = 1001.5
$ws.Range("L3").Value = 1001.5
$ws.Range("N3").Value = -1227.5
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 3924.5
$ws.Range("I62").Value = 2932.8
$ws.Range("J62").Value = 4632.857
$ws.Range("K62").Value = 2932.8
$ws.Range("L62").Value = 4632.857
$ws.Range("M62").Value = -2308.8
$ws.Range("N62").Value = -5880.857
# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 3924.5
$ws.Range("I65").Value = 2932.8
$ws.Range("J65").Value = 4632.857
$ws.Range("K65").Value = 14664
$ws.Range("L65").Value = 23164.285
$ws.Range("M65").Value = -11544
$ws.Range("N65").Value = -29404.285
# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water / Boiled Egg
$ws.Range("H4").Value = 163.66667
$ws.Range("I4").Value = 163.66667
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 491.00001
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -379.00001
$ws.Range("N4").ClearContents()
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 728.9388
$ws.Range("I131").Value = 385.94736
$ws.Range("J131").Value = 811.43036
$ws.Range("K131").Value = 1157.84208
$ws.Range("L131").Value = 2434.29108
$ws.Range("M131").Value = 3882.15792
$ws.Range("N131").Value = -12514.29108
# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 17: Point of Honor / Amateur's Needle
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 4234.273
$ws.Range("I70").Value = 4183.5
$ws.Range("J70").Value = 4369.6665
$ws.Range("K70").Value = 4183.5
$ws.Range("L70").Value = 4369.6665
$ws.Range("M70").Value = -3913.5
$ws.Range("N70").Value = -4909.6665
# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 4234.273
$ws.Range("I73").Value = 4183.5
$ws.Range("J73").Value = 4369.6665
$ws.Range("K73").Value = 4183.5
$ws.Range("L73").Value = 4369.6665
$ws.Range("M73").Value = -3247.5
$ws.Range("N73").Value = -6241.6665
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1949.875
$ws.Range("I102").Value = 1850
$ws.Range("J102").Value = 2049.75
$ws.Range("K102").Value = 1850
$ws.Range("L102").Value = 2049.75
$ws.Range("M102").Value = -228
$ws.Range("N102").Value = -5293.75
# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head / Leather Calot
$ws.Range("H2").Value = 252500
$ws.Range("I2").Value = 500000
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 500000
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -499888
$ws.Range("N2").Value = -5224
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 6649.3213
$ws.Range("I122").Value = 8722.706
$ws.Range("J122").Value = 3445
$ws.Range("K122").Value = 26168.118
$ws.Range("L122").Value = 10335
$ws.Range("M122").Value = -23718.118
$ws.Range("N122").Value = -15235
# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 9: A Taste for Dalmaticae / Amateur's Dalmatica
$ws.Range("H9").Value = 56671.332
$ws.Range("I9").Value = 10000
$ws.Range("J9").Value = 80007
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 80007
$ws.Range("M9").Value = -9860
$ws.Range("N9").Value = -80287
